$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill row 2 with sample import data
$ws.Range("A2").Value = 323230002
$ws.Range("B2").Value = 123
$ws.Range("C2").Value = "2023-11-01 07:50:00"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 111
$ws.Range("F2").Value = "2023-11-01 07:50:00"
$ws.Range("G2").Value = "2023-11-01 05:05:00"
$ws.Range("H2").Value = "2023-11-01 07:50:00"

$ws.Range("C2:H2").NumberFormat = "@"

# Re-apply explicit width to columns E:G (keep same width, make it "custom")
$ws.Range("E1:G1").ColumnWidth = 27.72

# Update selection/view state
$ws.Range("H24").Select()
